$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "SearchUser_Result2"
$ws.Range("B24").Value = "//td[normalize-space()='Admin']"
$ws.Range("C24").Value = "By.xpath"
